$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row: "_old" -> "_FV2304" suffix, "_new" -> "_FV2310" suffix.
$oldHeaders = @(
    "Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old",
    "Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old",
    "diff",
    "Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new",
    "Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new"
)
$newHeaders = @(
    "Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304",
    "Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310",
    "Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

# 2) Freeze the header row (split/freeze pane at row 2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Turn the used range into a real Excel table (ListObject) named Table1.
$tableRange = $ws.Range("A1:U55")
$table = $ws.ListObjects.Add(1, $tableRange, $false, 1, [System.Type]::Missing)
$table.Name = "Table1"
